$d = $word.ActiveDocument

$pairs = @(
    @("946×6=", "575×8="),
    @("532×5=", "254×3="),
    @("726×5=", "353×4="),
    @("105×2=", "870×8="),
    @("537×5=", "175×6="),
    @("161×3=", "593×6="),
    @("905×4=", "307×2="),
    @("910×4=", "653×7="),
    @("576×8=", "862×8="),
    @("769×8=", "579×4="),
    @("764×3=", "603×3="),
    @("618×7=", "506×6="),
    @("674×3=", "550×7="),
    @("811×9=", "845×2="),
    @("575×7=", "456×4="),
    @("407×4=", "780×2="),
    @("481×4=", "444×9="),
    @("482×5=", "333×3="),
    @("778×4=", "621×8="),
    @("778×8=", "180×4="),
    @("455×2=", "759×5="),
    @("783×5=", "392×5="),
    @("303×5=", "947×8="),
    @("139×2=", "487×2="),
    @("631×9=", "863×6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done replacing $($pairs.Count) cells"
